$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row after the "chemical_recycling_gasification" row (row 9),
# shifting all following rows down by one.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row with the new "chemical_recycling_pyrolysis" parameter.
$ws.Range("A10").Value = "chemical_recycling_pyrolysis"
$ws.Range("B10").Value = $true
